# "Works with multiple elements"
#
# The IFC_sign sheet lists one row per recognised element (ModelName /
# ElementId / GlobalId). Previously it only ever had a single data row;
# this update refreshes that row with a new element and appends three
# more rows so the sheet demonstrates handling several elements at once.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFC_sign")

# Helper: write a value that must be stored as literal text even though it
# looks numeric (e.g. "73368"). A leading apostrophe forces Excel to treat
# it as text instead of coercing it to a number; we then restore the
# cell's normal (non "quote-prefixed") look by re-applying the formatting
# already used by the sheet's other text cells.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
}

$globalId = "1hS`$VweOwvHwyp3_HHdqTH"

# Row 2 already existed - refresh it with the new element id/guid.
Set-TextValue $ws.Range("B2") "73366"
$ws.Range("C2").Value = $globalId

# Rows 3-5 are the additional elements.
$ws.Range("A3").Value = "Default"
Set-TextValue $ws.Range("B3") "73368"
$ws.Range("C3").Value = $globalId

$ws.Range("A4").Value = "Default"
Set-TextValue $ws.Range("B4") "73369"
$ws.Range("C4").Value = $globalId

$ws.Range("A5").Value = "Default"
Set-TextValue $ws.Range("B5") "73370"
$ws.Range("C5").Value = $globalId

# The apostrophe-prefix trick above marks B2:B5 with a "quote prefix"
# style; copy the plain formatting already used elsewhere on the sheet
# back onto those cells so they keep looking like ordinary text cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2:B5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Column C needs to be a little wider to comfortably fit the GlobalId text.
$ws.Columns.Item(3).ColumnWidth = 26.8
